$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.352.26"
$ws.Range("E2").Value = "  +8.79%  "
$ws.Range("D3").Value = "1.602.52"
$ws.Range("E3").Value = "  +8.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9977"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.96"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3717"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3387"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.34"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07058"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.85"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +9.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.941"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.666"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001085"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9964"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "1.597.25"
$ws.Range("E18").Value = "  +8.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06640"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +11.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +12.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +10.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.91%  "
$ws.Range("D24").Value = "22.314.63"
$ws.Range("E24").Value = "  +8.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.505"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +17.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +14.45%  "
$ws.Range("D29").Value = "1.776.49"
$ws.Range("E29").Value = "  +8.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.165"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.026"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +20.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9495"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +17.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08259"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("E35").Value = "  +5.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.308"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.666"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06158"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.242"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02220"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2028"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9959"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5952"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +12.41%  "
$ws.Range("E45").Value = "  +8.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.668"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5727"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.975"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +10.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06829"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.39%  "
